$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-15 from 2023-09-21 (45190) to 2023-09-23 (45192)
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
